$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price observation needs to be recorded between the existing
# row 71 and what was row 72, so insert a fresh row at position 72 which
# pushes all the following rows (old 72..167) down to (73..168).
$ws.Rows("72:72").Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Cells.Item(72, 1).Value2  = 1
$ws.Cells.Item(72, 2).Value2  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(72, 3).Value2  = "Arica y Parinacota"
$ws.Cells.Item(72, 4).Value2  = 45079
$ws.Cells.Item(72, 5).Value2  = 15
$ws.Cells.Item(72, 6).Value2  = 100112042
$ws.Cells.Item(72, 7).Value2  = "Locoto"
$ws.Cells.Item(72, 8).Value2  = "Sin especificar"
$ws.Cells.Item(72, 9).Value2  = "Segunda"
$ws.Cells.Item(72, 10).Value2 = 120
$ws.Cells.Item(72, 11).Value2 = 64000
$ws.Cells.Item(72, 12).Value2 = 65000
$ws.Cells.Item(72, 13).Value2 = 64500
$ws.Cells.Item(72, 14).Value2 = "$/caja 20 kilos"
$ws.Cells.Item(72, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(72, 16).Value2 = 3225
$ws.Cells.Item(72, 17).Value2 = 20
$ws.Cells.Item(72, 18).Value2 = "Hortaliza"
